$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.629.19"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "2.714.47"
$ws.Range("E3").Value = "  -6.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'503.75"
$ws.Range("E5").Value = "  -4.48%  "
$ws.Range("D6").Value = "'141.62"
$ws.Range("E6").Value = "  -0.56%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E8").Value = "  -3.97%  "
$ws.Range("D9").Value = "2.722.95"
$ws.Range("E9").Value = "  -6.42%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").Value = "'6.07"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("D12").Value = "'0.350"
$ws.Range("E12").Value = "  -2.46%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "3.189.91"
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("D15").Value = "58.650.56"
$ws.Range("E15").Value = "  -3.11%  "
$ws.Range("D16").Value = "'21.76"
$ws.Range("E16").Value = "  -3.44%  "
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").Value = "2.719.22"
$ws.Range("E18").Value = "  -6.48%  "
$ws.Range("E19").Value = "  -5.14%  "
$ws.Range("D20").Value = "'10.99"
$ws.Range("E20").Value = "  -5.57%  "
$ws.Range("D21").Value = "'343.02"
$ws.Range("E21").Value = "  -5.59%  "
$ws.Range("D22").Value = "'6.26"
$ws.Range("E22").Value = "  -4.58%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "'5.63"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").Value = "'63.20"
$ws.Range("E25").Value = "  -1.39%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.174"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").Value = "'0.428"
$ws.Range("E27").Value = "  -5.08%  "
$ws.Range("D29").Value = "'7.51"
$ws.Range("E29").Value = "  -3.21%  "
$ws.Range("D30").Value = "0.0₃0833"
$ws.Range("E30").Value = "  -2.84%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'19.27"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("E33").Value = "  -4.02%  "
$ws.Range("D34").Value = "'150.79"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").Value = "'4.22"
$ws.Range("E35").Value = "  -2.86%  "
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("D37").Value = "'0.948"
$ws.Range("E37").Value = "  -5.27%  "
$ws.Range("D38").Value = "'1.13"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  -6.90%  "
$ws.Range("D42").Value = "2.186.83"
$ws.Range("E42").Value = "  -5.89%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "'0.598"
$ws.Range("E45").Value = "  -7.05%  "
$ws.Range("D46").Value = "'19.07"
$ws.Range("E46").Value = "  -7.22%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.38"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'4.75"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("E49").Value = "  -2.60%  "
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").Value = "'18.16"
$ws.Range("E51").Value = "  -1.43%  "
